# Fruta / hortaliza, semanal
# Insert a new weekly record as row 7 (pushing the existing rows 7-68 down to 8-69).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(7).Insert()

$ws.Range("A7").Value = 1
$ws.Range("B7").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C7").Value = "Arica y Parinacota"
$ws.Range("D7").Value = 44959
$ws.Range("E7").Value = 15
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100103
$ws.Range("H7").Value = "Frutos de hueso (carozo)"
$ws.Range("I7").Value = 100103004
$ws.Range("J7").Value = "Durazno"
$ws.Range("K7").Value = "Springcrest"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 550
$ws.Range("N7").Value = 10000
$ws.Range("O7").Value = 11000
$ws.Range("P7").Value = 10364
$ws.Range("Q7").Value = "$/bandeja 10 kilos granel"
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 1036
$ws.Range("T7").Value = 10
